$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Submitted", $true, $false, $false, $false, $false, $true, 1, $false, "Reject and Resubmit at Review of Economic Studies", 2)
